$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate() | Out-Null

# The "Cases" query in B2 (CasesTab row) is being corrected: the trailing
# "Cohort" output column (and its value) is removed from the Cypher query
# text, while the rest of the query is unchanged.
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
   WHERE f.file_type IN ["Index File"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Restore the active selection on the sheet to B2 (matches the saved view
# state after the edit).
$ws.Range("B2").Select() | Out-Null
